$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fecha baja" header (column G) is actually meant to be "Fecha activacion";
# the real "Fecha baja" column already lives in H (it was a duplicated label).
$ws.Range("G1").Value = "Fecha activacion"

# Widen column G so the longer "Fecha activacion" label fits (~13.52 chars).
$ws.Columns.Item(7).ColumnWidth = 12.6

# The id columns L:P (Activo grupo id .. Metodo depreciacion id) should be
# stored/displayed as text, not as a general number.
$ws.Range("L2:P2").NumberFormat = "@"

# Move the active selection to F6.
$ws.Range("F6").Select()
